# Update cryptocurrency price values in column D on the active sheet.
# These cells hold their values as text (string) rather than numbers, so we
# temporarily mark the cell as Text-formatted before assigning the new
# value (otherwise a numeric-looking string like "267.72" gets auto-coerced
# to a Number by COM), then restore the cell's original ("Normal") style so
# no stray formatting is left behind.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = "267.72"
    3  = "21.55"
    4  = "6.192"
    5  = "0.06169"
    6  = "3.565"
    7  = "6.550"
    8  = "1.363"
    10 = "0.01344"
    12 = "0.08161"
    13 = "0.03320"
    14 = "0.03176"
    15 = "0.09268"
    16 = "3.747"
    17 = "0.001647"
    19 = "0.006338"
    20 = "0.005766"
    22 = "0.0001498"
    23 = "3.738"
    24 = "2.317"
    25 = "0.3302"
    28 = "0.0001615"
    40 = "0.04660"
    41 = "0.006976"
    42 = "0.1134"
    43 = "0.003655"
    44 = "0.01154"
    45 = "0.00005935"
    46 = "0.0009885"
    48 = "0.7808"
    49 = "0.002438"
    50 = "0.00001897"
    51 = "0.01238"
}

foreach ($row in $updates.Keys) {
    $cell = $ws.Cells.Item($row, 4)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$row]
    $cell.Style = "Normal"
}
